# Apply "new arrangement" edit: add SL NO / CONTEMP / IS MAIN columns to both
# sheets, number the acharyan_captions rows into the new CONTEMP column,
# group the Founders_Early_Acharyas rows into contemporaneous clusters (with
# an "M" marker for the main acharya in a cluster), swap two mis-ordered
# names, and drop the now-unused trailing blank rows on sheet 2.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # acharyan_captions
$ws2 = $wb.Worksheets.Item(2)   # Founders_Early_Acharyas

$xlCenter = -4108

# ---------------------------------------------------------------------------
# Seed the shared-string table in the exact order the target workbook uses
# (M=56, SL NO=57, CONTEMP=58, IS MAIN=59) by writing the "M" marker first.
# ---------------------------------------------------------------------------
$ws2.Cells.Item(10, 4).Value = "M"
$ws2.Cells.Item(18, 4).Value = "M"
$ws2.Range("D10").HorizontalAlignment = $xlCenter
$ws2.Range("D10").VerticalAlignment = $xlCenter
$ws2.Range("D18").HorizontalAlignment = $xlCenter
$ws2.Range("D18").VerticalAlignment = $xlCenter

# ---------------------------------------------------------------------------
# Sheet 1: acharyan_captions
# ---------------------------------------------------------------------------
$ws1.Cells.Item(1, 1).Value = "SL NO"

$ws1.Cells.Item(1, 3).Value = "CONTEMP"
$ws1.Range("C1").HorizontalAlignment = $xlCenter
$ws1.Range("C1").VerticalAlignment = $xlCenter

$ws1.Cells.Item(1, 4).Value = "IS MAIN"
$ws1.Range("D1").HorizontalAlignment = $xlCenter
$ws1.Range("D1").VerticalAlignment = $xlCenter

for ($row = 2; $row -le 37; $row++) {
    $ws1.Cells.Item($row, 3).Value = $row - 1
}

$ws1.Range("C2:C37").Select() | Out-Null
$ws1.Application.ActiveWindow.RangeSelection | Out-Null

# ---------------------------------------------------------------------------
# Sheet 2: Founders_Early_Acharyas
# ---------------------------------------------------------------------------
$ws2.Cells.Item(1, 1).Value = "SL NO"
$ws2.Cells.Item(1, 3).Value = "CONTEMP"
$ws2.Cells.Item(1, 4).Value = "IS MAIN"
$ws2.Range("A1").HorizontalAlignment = $xlCenter
$ws2.Range("A1").VerticalAlignment = $xlCenter
$ws2.Range("C1").HorizontalAlignment = $xlCenter
$ws2.Range("C1").VerticalAlignment = $xlCenter
$ws2.Range("D1").HorizontalAlignment = $xlCenter
$ws2.Range("D1").VerticalAlignment = $xlCenter

$contemp = @(0,1,2,3,4,5,6,7,8,9,9,10,10,11,11,12,13,14)
for ($i = 0; $i -lt $contemp.Length; $i++) {
    $row = $i + 2
    $ws2.Cells.Item($row, 3).Value = $contemp[$i]
    $ws2.Range("C$row").HorizontalAlignment = $xlCenter
    $ws2.Range("C$row").VerticalAlignment = $xlCenter
}

# Two names were swapped into the correct contemporaneous rows.
$ws2.Cells.Item(14, 2).Value = " Engan Āḷvān"
$ws2.Cells.Item(15, 2).Value = " Sri Kidambi Raṅgarāja Āchārya"

# The sheet used to carry a pile of trailing, completely blank placeholder
# rows (20-35) below the real data; the new layout drops them.
$ws2.Rows("20:35").Delete() | Out-Null

# Column widths/centring for the two new helper columns.
$ws2.Range("C1:D19").EntireColumn.ColumnWidth = 8.04
$ws2.Range("C1:D19").EntireColumn.HorizontalAlignment = $xlCenter
$ws2.Range("C1:D19").EntireColumn.VerticalAlignment = $xlCenter

$ws2.Range("A1:D1").Select() | Out-Null

Write-Output "edit applied"
